# Household member survey workbook edit:
#  - Remove the "for {{instance_name}}" suffix from the age note template
#    (now that instance_name will be a setting instead of a data value).
#  - Add a new 'instance_name' setting row pointing at 'member_name'.
#  - Update view/selection state: settings tab becomes the active tab,
#    survey tab is no longer the selected one, and selections/row height
#    are adjusted to match.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey sheet: update the note text in D6, grow the row to fit, and
#     move the selection off of E7 onto D7 ---
$survey.Range("D6").Value = "{{member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}}"
$survey.Rows.Item(6).RowHeight = 62.5

# --- settings sheet: add a new 'instance_name' -> 'member_name' row ---
$settings.Range("A6").Value = "instance_name"
$settings.Range("B6").Value = "member_name"
$settings.Rows.Item(6).RowHeight = 12.75

# --- selections / active sheet: settings becomes the active tab, survey's
#     selection moves to D7, settings' selection moves to B18 ---
$survey.Range("D7").Select()
$settings.Activate()
$settings.Range("B18").Select()
